$d = $word.ActiveDocument
$normal = $d.Styles("Normal")
# wdFrench = 1036, wdEnglishUS = 1033, wdArabic = 1025
$normal.Font.LanguageID = 1036
$normal.Font.LanguageIDFarEast = 1033
$normal.Font.LanguageIDBi = 1025
Write-Output $normal.Font.LanguageID
Write-Output $normal.Font.LanguageIDFarEast
Write-Output $normal.Font.LanguageIDBi
